$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(87, 8).Value = 44677
$ws.Cells.Item(87, 10).Value = 44677
$ws.Cells.Item(87, 12).Value = 44677
$ws.Cells.Item(87, 14).Value = -47173
$ws.Cells.Item(90, 8).Value = 44677
$ws.Cells.Item(90, 10).Value = 44677
$ws.Cells.Item(90, 12).Value = 134031
$ws.Cells.Item(90, 14).Value = -146511
$ws.Cells.Item(108, 8).Value = 29612.5
$ws.Cells.Item(108, 10).Value = 29612.5
$ws.Cells.Item(108, 12).Value = 29612.5
$ws.Cells.Item(108, 14).Value = -37292.5
$ws.Cells.Item(114, 8).Value = 36698
$ws.Cells.Item(114, 10).Value = 36698
$ws.Cells.Item(114, 12).Value = 36698
$ws.Cells.Item(114, 14).Value = -45376
$ws.Cells.Item(130, 8).Value = 43623.2
$ws.Cells.Item(130, 10).Value = 43623.2
$ws.Cells.Item(130, 12).Value = 43623.2
$ws.Cells.Item(130, 14).Value = -53663.2

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(34, 8).Value = 2500
$ws.Cells.Item(34, 9).Value = 2500
$ws.Cells.Item(34, 10).Value = 0
$ws.Cells.Item(34, 11).Value = 2500
$ws.Cells.Item(34, 12).Value = 0
$ws.Cells.Item(34, 13).Value = -2229
$ws.Cells.Item(34, 14).ClearContents()
$ws.Cells.Item(61, 8).Value = 2466.1296
$ws.Cells.Item(61, 9).Value = 1406.6757
$ws.Cells.Item(61, 10).Value = 4772
$ws.Cells.Item(61, 11).Value = 1406.6757
$ws.Cells.Item(61, 12).Value = 4772
$ws.Cells.Item(61, 13).Value = -1194.6757
$ws.Cells.Item(61, 14).Value = -5196
$ws.Cells.Item(110, 8).Value = 1405.5555
$ws.Cells.Item(110, 9).Value = 1385.0869
$ws.Cells.Item(110, 11).Value = 1385.0869
$ws.Cells.Item(110, 13).Value = 659.9131
$ws.Cells.Item(119, 8).Value = 50698
$ws.Cells.Item(119, 10).Value = 50698
$ws.Cells.Item(119, 12).Value = 50698
$ws.Cells.Item(119, 14).Value = -60374
$ws.Cells.Item(131, 8).Value = 49992
$ws.Cells.Item(131, 10).Value = 49992
$ws.Cells.Item(131, 12).Value = 49992
$ws.Cells.Item(131, 14).Value = -60072
$ws.Cells.Item(133, 8).Value = 46113.875
$ws.Cells.Item(133, 10).Value = 46113.875
$ws.Cells.Item(133, 12).Value = 46113.875
$ws.Cells.Item(133, 14).Value = -51173.875
$ws.Cells.Item(136, 8).Value = 2466.1296
$ws.Cells.Item(136, 9).Value = 1406.6757
$ws.Cells.Item(136, 10).Value = 4772
$ws.Cells.Item(136, 11).Value = 4220.0271
$ws.Cells.Item(136, 12).Value = 14316
$ws.Cells.Item(136, 13).Value = -1670.0271
$ws.Cells.Item(136, 14).Value = -19416
$ws.Cells.Item(137, 8).Value = 31264.285
$ws.Cells.Item(137, 10).Value = 31264.285
$ws.Cells.Item(137, 12).Value = 31264.285
$ws.Cells.Item(137, 14).Value = -41464.285

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(122, 8).Value = 37138.4
$ws.Cells.Item(122, 10).Value = 37138.4
$ws.Cells.Item(122, 12).Value = 37138.4
$ws.Cells.Item(122, 14).Value = -46938.4

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(16, 8).Value = 1548.7
$ws.Cells.Item(16, 9).Value = 892.2
$ws.Cells.Item(16, 10).Value = 2205.2
$ws.Cells.Item(16, 11).Value = 892.2
$ws.Cells.Item(16, 12).Value = 2205.2
$ws.Cells.Item(16, 13).Value = -605.2
$ws.Cells.Item(16, 14).Value = -2779.2
$ws.Cells.Item(31, 8).Value = 264545.38
$ws.Cells.Item(31, 9).Value = 4038
$ws.Cells.Item(31, 10).Value = 291127.78
$ws.Cells.Item(31, 11).Value = 4038
$ws.Cells.Item(31, 12).Value = 291127.78
$ws.Cells.Item(31, 13).Value = -3743
$ws.Cells.Item(31, 14).Value = -291717.78
$ws.Cells.Item(34, 8).Value = 264545.38
$ws.Cells.Item(34, 9).Value = 4038
$ws.Cells.Item(34, 10).Value = 291127.78
$ws.Cells.Item(34, 11).Value = 4038
$ws.Cells.Item(34, 12).Value = 291127.78
$ws.Cells.Item(34, 13).Value = -3836
$ws.Cells.Item(34, 14).Value = -291531.78
$ws.Cells.Item(99, 8).Value = 1861.6
$ws.Cells.Item(99, 9).Value = 1544.8
$ws.Cells.Item(99, 11).Value = 1544.8
$ws.Cells.Item(99, 13).Value = -46.79999999999995
$ws.Cells.Item(111, 8).Value = 46997.332
$ws.Cells.Item(111, 10).Value = 46997.332
$ws.Cells.Item(111, 12).Value = 46997.332
$ws.Cells.Item(111, 14).Value = -55177.332
$ws.Cells.Item(112, 8).Value = 32235.143
$ws.Cells.Item(112, 10).Value = 32235.143
$ws.Cells.Item(112, 12).Value = 32235.143
$ws.Cells.Item(112, 14).Value = -35189.143
$ws.Cells.Item(113, 8).Value = 1548.7
$ws.Cells.Item(113, 9).Value = 892.2
$ws.Cells.Item(113, 10).Value = 2205.2
$ws.Cells.Item(113, 11).Value = 892.2
$ws.Cells.Item(113, 12).Value = 2205.2
$ws.Cells.Item(113, 13).Value = 1277.8
$ws.Cells.Item(113, 14).Value = -6545.2
$ws.Cells.Item(126, 8).Value = 1861.6
$ws.Cells.Item(126, 9).Value = 1544.8
$ws.Cells.Item(126, 11).Value = 4634.4
$ws.Cells.Item(126, 13).Value = -2164.4
$ws.Cells.Item(134, 8).Value = 584489.3
$ws.Cells.Item(134, 9).Value = 869.2222
$ws.Cells.Item(134, 10).Value = 2335349.5
$ws.Cells.Item(134, 11).Value = 2607.6666
$ws.Cells.Item(134, 12).Value = 7006048.5
$ws.Cells.Item(134, 13).Value = -72.66660000000002
$ws.Cells.Item(134, 14).Value = -7011118.5
$ws.Cells.Item(137, 8).Value = 18400
$ws.Cells.Item(137, 10).Value = 18400
$ws.Cells.Item(137, 12).Value = 18400
$ws.Cells.Item(137, 14).Value = -28600

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(113, 8).Value = 3755.5
$ws.Cells.Item(113, 9).Value = 5599.05
$ws.Cells.Item(113, 10).Value = 682.9167
$ws.Cells.Item(113, 11).Value = 16797.15
$ws.Cells.Item(113, 12).Value = 2048.7501
$ws.Cells.Item(113, 13).Value = -14627.15
$ws.Cells.Item(113, 14).Value = -6388.7501

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(28, 8).Value = 34997.5
$ws.Cells.Item(28, 10).Value = 39996.668
$ws.Cells.Item(28, 12).Value = 39996.668
$ws.Cells.Item(28, 14).Value = -40380.668
$ws.Cells.Item(33, 8).Value = 29000
$ws.Cells.Item(33, 9).Value = 0
$ws.Cells.Item(33, 10).Value = 29000
$ws.Cells.Item(33, 11).Value = 0
$ws.Cells.Item(33, 12).Value = 29000
$ws.Cells.Item(33, 13).ClearContents()
$ws.Cells.Item(33, 14).Value = -29504
$ws.Cells.Item(97, 8).Value = 5306.154
$ws.Cells.Item(97, 9).Value = 1470.625
$ws.Cells.Item(97, 11).Value = 1470.625
$ws.Cells.Item(97, 13).Value = -974.625
$ws.Cells.Item(113, 8).Value = 1701.963
$ws.Cells.Item(113, 9).Value = 1700.0714
$ws.Cells.Item(113, 10).Value = 1704
$ws.Cells.Item(113, 11).Value = 1700.0714
$ws.Cells.Item(113, 12).Value = 1704
$ws.Cells.Item(113, 13).Value = 469.9286
$ws.Cells.Item(113, 14).Value = -6044
$ws.Cells.Item(130, 8).Value = 49035.332
$ws.Cells.Item(130, 10).Value = 49035.332
$ws.Cells.Item(130, 12).Value = 49035.332
$ws.Cells.Item(130, 14).Value = -59075.332
$ws.Cells.Item(138, 8).Value = 38742.855
$ws.Cells.Item(138, 10).Value = 38742.855
$ws.Cells.Item(138, 12).Value = 38742.855
$ws.Cells.Item(138, 14).Value = -49022.855

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(99, 8).Value = 33166.668
$ws.Cells.Item(99, 9).Value = 25000
$ws.Cells.Item(99, 11).Value = 25000
$ws.Cells.Item(99, 13).Value = -22005
$ws.Cells.Item(125, 8).Value = 43230.668
$ws.Cells.Item(125, 10).Value = 43230.668
$ws.Cells.Item(125, 12).Value = 43230.668
$ws.Cells.Item(125, 14).Value = -53070.668
$ws.Cells.Item(135, 8).Value = 35899
$ws.Cells.Item(135, 10).Value = 35899
$ws.Cells.Item(135, 12).Value = 35899
$ws.Cells.Item(135, 14).Value = -46039
$ws.Cells.Item(137, 8).Value = 36980
$ws.Cells.Item(137, 10).Value = 36980
$ws.Cells.Item(137, 12).Value = 36980
$ws.Cells.Item(137, 14).Value = -47180
$ws.Cells.Item(139, 8).Value = 46562.25
$ws.Cells.Item(139, 10).Value = 46562.25
$ws.Cells.Item(139, 12).Value = 46562.25
$ws.Cells.Item(139, 14).Value = -56842.25

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(115, 8).Value = 37377
$ws.Cells.Item(115, 10).Value = 37377
$ws.Cells.Item(115, 12).Value = 37377
$ws.Cells.Item(115, 14).Value = -40511
$ws.Cells.Item(116, 8).Value = 49680
$ws.Cells.Item(116, 10).Value = 49680
$ws.Cells.Item(116, 12).Value = 49680
$ws.Cells.Item(116, 14).Value = -58858
$ws.Cells.Item(118, 8).Value = 34592
$ws.Cells.Item(118, 10).Value = 44388
$ws.Cells.Item(118, 12).Value = 44388
$ws.Cells.Item(118, 14).Value = -47702
$ws.Cells.Item(121, 8).Value = 43262
$ws.Cells.Item(121, 10).Value = 43262
$ws.Cells.Item(121, 12).Value = 43262
$ws.Cells.Item(121, 14).Value = -46756
$ws.Cells.Item(125, 8).Value = 35474
$ws.Cells.Item(125, 10).Value = 35474
$ws.Cells.Item(125, 12).Value = 35474
$ws.Cells.Item(125, 14).Value = -45314
$ws.Cells.Item(139, 8).Value = 44466.668
$ws.Cells.Item(139, 10).Value = 44466.668
$ws.Cells.Item(139, 12).Value = 44466.668
$ws.Cells.Item(139, 14).Value = -54746.668
